# ConfigCategorie.xlsx edit
# Commit message: "Changement du redimenssionement de la page gmail à l'ouverture"
#
# The EMPLOI row's value list (cell B4, column "Value") gains a new
# entry "mission-locale" at the end, which in turn makes Excel's
# auto row-height (wrap text is on) grow from 58 to 72.5 points once
# the cell is saved, since the wrapped text now needs more lines.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Categories")

# Update the "Value" cell for the EMPLOI category, appending the
# new keyword "mission-locale" to the existing semicolon separated list.
$ws.Range("B4").Value = "Talent.com;HelloWork;Jungle;Linkedin;capgemini;emploi;meteojob;mission-locale"

# The longer text now wraps onto more lines inside the cell, so the
# row grows taller (58 -> 72.5 points) to keep showing the full value.
$ws.Rows.Item(4).RowHeight = 72.5
